$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2314.0408
$ws.Range("I15").Value = 2314.0408
$ws.Range("K15").Value = 6942.1224
$ws.Range("M15").Value = -6773.1224

$ws.Range("H74").Value = 3939.2307
$ws.Range("I74").Value = 3867
$ws.Range("J74").Value = 4001.1428
$ws.Range("K74").Value = 3867
$ws.Range("L74").Value = 4001.1428
$ws.Range("M74").Value = -2931
$ws.Range("N74").Value = -5873.1428

$ws.Range("H77").Value = 3939.2307
$ws.Range("I77").Value = 3867
$ws.Range("J77").Value = 4001.1428
$ws.Range("K77").Value = 19335
$ws.Range("L77").Value = 20005.714
$ws.Range("M77").Value = -14655
$ws.Range("N77").Value = -29365.714

$ws.Range("H112").Value = 27779340
$ws.Range("J112").Value = 1810.1
$ws.Range("L112").Value = 5430.299999999999
$ws.Range("N112").Value = -7646.299999999999

$ws.Range("H126").Value = 49200
$ws.Range("J126").Value = 49200
$ws.Range("L126").Value = 49200
$ws.Range("N126").Value = -59080

$ws.Range("H129").Value = 805.6609999999999
$ws.Range("I129").Value = 403
$ws.Range("J129").Value = 968.6429000000001
$ws.Range("K129").Value = 1209
$ws.Range("L129").Value = 2905.9287
$ws.Range("M129").Value = 3791
$ws.Range("N129").Value = -12905.9287

$ws.Range("H132").Value = 5448205.5
$ws.Range("I132").Value = 6785
$ws.Range("J132").Value = 9801342
$ws.Range("K132").Value = 20355
$ws.Range("L132").Value = 29404026
$ws.Range("M132").Value = -17825
$ws.Range("N132").Value = -29409086

$ws.Range("H138").Value = 1014084.5
$ws.Range("I138").Value = 1383.08
$ws.Range("J138").Value = 1951771
$ws.Range("K138").Value = 4149.24
$ws.Range("L138").Value = 5855313
$ws.Range("M138").Value = 990.7600000000002
$ws.Range("N138").Value = -5865593

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1941.4
$ws.Range("I45").Value = 3400
$ws.Range("K45").Value = 3400
$ws.Range("M45").Value = -3023

$ws.Range("H61").Value = 143286320
$ws.Range("I61").Value = 200200850
$ws.Range("J61").Value = 999999.5
$ws.Range("K61").Value = 200200850
$ws.Range("L61").Value = 999999.5
$ws.Range("M61").Value = -200200638
$ws.Range("N61").Value = -1000423.5

$ws.Range("H132").Value = 66059.30499999999
$ws.Range("I132").Value = 51882.95
$ws.Range("J132").Value = 87869.08
$ws.Range("K132").Value = 155648.85
$ws.Range("L132").Value = 263607.24
$ws.Range("M132").Value = -153118.85
$ws.Range("N132").Value = -268667.24

$ws.Range("H136").Value = 143286320
$ws.Range("I136").Value = 200200850
$ws.Range("J136").Value = 999999.5
$ws.Range("K136").Value = 600602550
$ws.Range("L136").Value = 2999998.5
$ws.Range("M136").Value = -600600000
$ws.Range("N136").Value = -3005098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5580.4
$ws.Range("I134").Value = 5725.5
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 17176.5
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -14641.5
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 71431336
$ws.Range("I58").Value = 250001000
$ws.Range("J58").Value = 3465.5
$ws.Range("K58").Value = 250001000
$ws.Range("L58").Value = 3465.5
$ws.Range("M58").Value = -250000797
$ws.Range("N58").Value = -3871.5

$ws.Range("H132").Value = 27483.615
$ws.Range("I132").Value = 1489.4117
$ws.Range("J132").Value = 47570.047
$ws.Range("K132").Value = 4468.2351
$ws.Range("L132").Value = 142710.141
$ws.Range("M132").Value = -1938.2351
$ws.Range("N132").Value = -147770.141

$ws.Range("H134").Value = 53732.76
$ws.Range("I134").Value = 3069.4
$ws.Range("J134").Value = 99790.37
$ws.Range("K134").Value = 9208.200000000001
$ws.Range("L134").Value = 299371.11
$ws.Range("M134").Value = -6673.200000000001
$ws.Range("N134").Value = -304441.11

$ws.Range("H136").Value = 71431336
$ws.Range("I136").Value = 250001000
$ws.Range("J136").Value = 3465.5
$ws.Range("K136").Value = 750003000
$ws.Range("L136").Value = 10396.5
$ws.Range("M136").Value = -750000450
$ws.Range("N136").Value = -15496.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 670.1875
$ws.Range("I5").Value = 548.2
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 1644.6
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -1532.6
$ws.Range("N5").Value = -7724

$ws.Range("H34").Value = 984.53845
$ws.Range("I34").Value = 433.33334
$ws.Range("J34").Value = 1149.9
$ws.Range("K34").Value = 1300.00002
$ws.Range("L34").Value = 3449.7
$ws.Range("M34").Value = -1216.00002
$ws.Range("N34").Value = -3617.7

$ws.Range("H39").Value = 799.3333
$ws.Range("J39").Value = 799.3333
$ws.Range("L39").Value = 2397.9999
$ws.Range("N39").Value = -2985.9999

$ws.Range("H55").Value = 2916.5833
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 3172.6365
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 9517.9095
$ws.Range("M55").Value = -123
$ws.Range("N55").Value = -9871.9095

$ws.Range("H122").Value = 18519330
$ws.Range("J122").Value = 41667892
$ws.Range("L122").Value = 375011028
$ws.Range("N122").Value = -375015928

$ws.Range("H135").Value = 670.1875
$ws.Range("I135").Value = 548.2
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 4933.8
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -2398.8
$ws.Range("N135").Value = -27570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1491.625
$ws.Range("I122").Value = 1490.4286
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4471.2858
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2021.2858
$ws.Range("N122").Value = -9400

$ws.Range("H132").Value = 224655.33
$ws.Range("I132").Value = 1000000
$ws.Range("J132").Value = 127737.25
$ws.Range("K132").Value = 3000000
$ws.Range("L132").Value = 383211.75
$ws.Range("M132").Value = -2997470
$ws.Range("N132").Value = -388271.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 933.3333
$ws.Range("I46").Value = 933.3333
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 933.3333
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -745.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 22991.143
$ws.Range("J92").Value = 22991.143
$ws.Range("L92").Value = 22991.143
$ws.Range("N92").Value = -27983.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N46").ClearContents()
